$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1)
# Rename Spanish accented municipality/state name particles to title case,
# fix "MonteMorelos" typo, and remove trailing footnote rows (723-727).

$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B20').Value = 'Amatenango De La Frontera'
$ws.Range('B25').Value = 'Comitán De Domínguez'
$ws.Range('B33').Value = 'Mazapa De Madero'
$ws.Range('B60').Value = 'Guadalupe Y Calvo'
$ws.Range('B62').Value = 'Hidalgo Del Parral'
$ws.Range('B88').Value = 'Villa De Álvarez'
$ws.Range('A90').Value = 'Ciudad De México'
$ws.Range('B94').Value = 'Cuajimalpa De Morelos'
$ws.Range('B117').Value = 'San Juan Del Río'
$ws.Range('B118').Value = 'San Luis Del Cordero'
$ws.Range('A124').Value = 'Estado De México'
$ws.Range('B124').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B126').Value = 'Almoloya De Juárez'
$ws.Range('B130').Value = 'Atizapán De Zaragoza'
$ws.Range('B133').Value = 'Chapa De Mota'
$ws.Range('B136').Value = 'Coacalco De Berriozábal'
$ws.Range('B140').Value = 'Ecatepec De Morelos'
$ws.Range('B146').Value = 'Ixtapan De La Sal'
$ws.Range('B147').Value = 'Ixtapan Del Oro'
$ws.Range('B151').Value = 'Naucalpan De Juárez'
$ws.Range('B157').Value = 'San Felipe Del Progreso'
$ws.Range('B160').Value = 'Tenango Del Valle'
$ws.Range('B163').Value = 'Tlalnepantla De Baz'
$ws.Range('B167').Value = 'Villa De Allende'
$ws.Range('B168').Value = 'Villa Del Carbón'
$ws.Range('B176').Value = 'Apaseo El Alto'
$ws.Range('B182').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B185').Value = 'Jaral Del Progreso'
$ws.Range('B190').Value = 'Purísima Del Rincón'
$ws.Range('B195').Value = 'San Francisco Del Rincón'
$ws.Range('B196').Value = 'San Luis De La Paz'
$ws.Range('B197').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B199').Value = 'Silao De La Victoria'
$ws.Range('B202').Value = 'Valle De Santiago'
$ws.Range('B207').Value = 'Acapulco De Juárez'
$ws.Range('B210').Value = 'Atoyac De Álvarez'
$ws.Range('B211').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B212').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B215').Value = 'Coyuca De Benítez'
$ws.Range('B216').Value = 'Coyuca De Catalán'
$ws.Range('B219').Value = 'Cutzamala De Pinzón'
$ws.Range('B222').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B223').Value = 'Iguala De La Independencia'
$ws.Range('B224').Value = 'Zihuatanejo De Azueta'
$ws.Range('B225').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B233').Value = 'Taxco De Alarcón'
$ws.Range('B234').Value = 'Técpan De Galeana'
$ws.Range('B236').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B239').Value = 'Tlapa De Comonfort'
$ws.Range('B245').Value = 'Agua Blanca De Iturbide'
$ws.Range('B248').Value = 'Atotonilco El Grande'
$ws.Range('B249').Value = 'Cuautepec De Hinojosa'
$ws.Range('B251').Value = 'Jacala De Ledezma'
$ws.Range('B253').Value = 'Mineral Del Monte'
$ws.Range('B254').Value = 'Nopala De Villagrán'
$ws.Range('B255').Value = 'Pachuca De Soto'
$ws.Range('B260').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B263').Value = 'Tula De Allende'
$ws.Range('B264').Value = 'Tulancingo De Bravo'
$ws.Range('B265').Value = 'Zacualtipán De Ángeles'
$ws.Range('B271').Value = 'Atotonilco El Alto'
$ws.Range('B273').Value = 'Autlán De Navarro'
$ws.Range('B280').Value = 'Concepción De Buenos Aires'
$ws.Range('B281').Value = 'Cuautitlán De García Barragán'
$ws.Range('B288').Value = 'Huejuquilla El Alto'
$ws.Range('B293').Value = 'Lagos De Moreno'
$ws.Range('B296').Value = 'Ojuelos De Jalisco'
$ws.Range('B300').Value = 'San Cristóbal De La Barranca'
$ws.Range('B301').Value = 'San Juan De Los Lagos'
$ws.Range('B303').Value = 'San Martín De Bolaños'
$ws.Range('B305').Value = 'San Miguel El Alto'
$ws.Range('B306').Value = 'San Sebastián Del Oeste'
$ws.Range('B308').Value = 'Talpa De Allende'
$ws.Range('B309').Value = 'Tamazula De Gordiano'
$ws.Range('B315').Value = 'Teocuitatlán De Corona'
$ws.Range('B316').Value = 'Tepatitlán De Morelos'
$ws.Range('B317').Value = 'Tizapán El Alto'
$ws.Range('B318').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B323').Value = 'Unión De San Antonio'
$ws.Range('B324').Value = 'Unión De Tula'
$ws.Range('B327').Value = 'Yahualica De González Gallo'
$ws.Range('B328').Value = 'Zacoalco De Torres'
$ws.Range('B331').Value = 'Zapotitlán De Vadillo'
$ws.Range('B332').Value = 'Zapotlán Del Rey'
$ws.Range('B333').Value = 'Zapotlán El Grande'
$ws.Range('B348').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B420').Value = 'Puente De Ixtla'
$ws.Range('B427').Value = 'Amatlán De Cañas'
$ws.Range('B428').Value = 'Bahía De Banderas'
$ws.Range('B432').Value = 'Ixtlán Del Río'
$ws.Range('B438').Value = 'Santa María Del Oro'
$ws.Range('B448').Value = 'Montemorelos'
$ws.Range('B450').Value = 'San Nicolás De Los Garza'
$ws.Range('B456').Value = 'Coicoyán De Las Flores'
$ws.Range('B458').Value = 'Guevea De Humboldt'
$ws.Range('B459').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B460').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B461').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B463').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B464').Value = 'Oaxaca De Juárez'
$ws.Range('B465').Value = 'Ocotlán De Morelos'
$ws.Range('B466').Value = 'Putla Villa De Guerrero'
$ws.Range('B483').Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range('B489').Value = 'Santa Inés De Zaragoza'
$ws.Range('B507').Value = 'Santo Domingo De Morelos'
$ws.Range('B511').Value = 'Teotitlán De Flores Magón'
$ws.Range('B512').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B513').Value = 'Tlacolula De Matamoros'
$ws.Range('B514').Value = 'Tlalixtac De Cabrera'
$ws.Range('B516').Value = 'Villa De Etla'
$ws.Range('B517').Value = 'Zapotitlán Del Río'
$ws.Range('B532').Value = 'Huitzilan De Serdán'
$ws.Range('B534').Value = 'Izúcar De Matamoros'
$ws.Range('B539').Value = 'Palmar De Bravo'
$ws.Range('B548').Value = 'Tepexi De Rodríguez'
$ws.Range('B559').Value = 'Cadereyta De Montes'
$ws.Range('B561').Value = 'Pinal De Amoles'
$ws.Range('B564').Value = 'San Juan Del Río'
$ws.Range('B571').Value = 'Mexquitic De Carmona'
$ws.Range('B577').Value = 'Santa María Del Río'
$ws.Range('B581').Value = 'Villa De Ramos'
$ws.Range('B582').Value = 'Villa De Reyes'
$ws.Range('B641').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B657').Value = 'Ixhuatlán De Madero'
$ws.Range('B658').Value = 'Juchique De Ferrer'
$ws.Range('B663').Value = 'Ozuluama De Mascareñas'
$ws.Range('B666').Value = 'Poza Rica De Hidalgo'
$ws.Range('B671').Value = 'Sayula De Alemán'
$ws.Range('B680').Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range('B690').Value = 'Concepción Del Oro'
$ws.Range('B697').Value = 'Jiménez Del Teul'
$ws.Range('B702').Value = 'Nochistlán De Mejía'
$ws.Range('B710').Value = 'Teúl De González Ortega'
$ws.Range('B711').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B715').Value = 'Villa De Cos'

# Remove the trailing footnote/metadata rows so the sheet dimension shrinks back to A1:D721
$ws.Range('A723:D727').Clear()
